$d = $word.ActiveDocument

# Update the date heading (unique text in the document, safe to Find/Replace)
$d.Content.Find.Execute("2024-01-04 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-01-05 Friday", 2) | Out-Null

# Update every equation cell in the 20x5 practice table directly by position
# (direct cell addressing avoids any ambiguity from repeated/substring text values)
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "19+35="
$t.Cell(1, 2).Range.Text = "35-27="
$t.Cell(1, 3).Range.Text = "71-15="
$t.Cell(1, 4).Range.Text = "91-89="
$t.Cell(1, 5).Range.Text = "38+46="
$t.Cell(2, 1).Range.Text = "58+36="
$t.Cell(2, 2).Range.Text = "55+8="
$t.Cell(2, 3).Range.Text = "86+8="
$t.Cell(2, 4).Range.Text = "44+48="
$t.Cell(2, 5).Range.Text = "9+75="
$t.Cell(3, 1).Range.Text = "25+49="
$t.Cell(3, 2).Range.Text = "16+67="
$t.Cell(3, 3).Range.Text = "37+8="
$t.Cell(3, 4).Range.Text = "17+44="
$t.Cell(3, 5).Range.Text = "39+55="
$t.Cell(4, 1).Range.Text = "74-47="
$t.Cell(4, 2).Range.Text = "3+9="
$t.Cell(4, 3).Range.Text = "36+15="
$t.Cell(4, 4).Range.Text = "40-3="
$t.Cell(4, 5).Range.Text = "38+47="
$t.Cell(5, 1).Range.Text = "83-48="
$t.Cell(5, 2).Range.Text = "4+59="
$t.Cell(5, 3).Range.Text = "15-7="
$t.Cell(5, 4).Range.Text = "7+69="
$t.Cell(5, 5).Range.Text = "92-54="
$t.Cell(6, 1).Range.Text = "20-6="
$t.Cell(6, 2).Range.Text = "9+49="
$t.Cell(6, 3).Range.Text = "21-13="
$t.Cell(6, 4).Range.Text = "92-15="
$t.Cell(6, 5).Range.Text = "28+36="
$t.Cell(7, 1).Range.Text = "70-4="
$t.Cell(7, 2).Range.Text = "59+8="
$t.Cell(7, 3).Range.Text = "54-9="
$t.Cell(7, 4).Range.Text = "16+8="
$t.Cell(7, 5).Range.Text = "65+18="
$t.Cell(8, 1).Range.Text = "91-13="
$t.Cell(8, 2).Range.Text = "79+14="
$t.Cell(8, 3).Range.Text = "56-49="
$t.Cell(8, 4).Range.Text = "81-29="
$t.Cell(8, 5).Range.Text = "5+39="
$t.Cell(9, 1).Range.Text = "74+19="
$t.Cell(9, 2).Range.Text = "59+19="
$t.Cell(9, 3).Range.Text = "82-7="
$t.Cell(9, 4).Range.Text = "38-29="
$t.Cell(9, 5).Range.Text = "56-39="
$t.Cell(10, 1).Range.Text = "3+48="
$t.Cell(10, 2).Range.Text = "61-2="
$t.Cell(10, 3).Range.Text = "64-45="
$t.Cell(10, 4).Range.Text = "63-27="
$t.Cell(10, 5).Range.Text = "44+49="
$t.Cell(11, 1).Range.Text = "36+16="
$t.Cell(11, 2).Range.Text = "55-17="
$t.Cell(11, 3).Range.Text = "27+66="
$t.Cell(11, 4).Range.Text = "57+9="
$t.Cell(11, 5).Range.Text = "71-54="
$t.Cell(12, 1).Range.Text = "65-16="
$t.Cell(12, 2).Range.Text = "5+66="
$t.Cell(12, 3).Range.Text = "90-24="
$t.Cell(12, 4).Range.Text = "97-9="
$t.Cell(12, 5).Range.Text = "6+6="
$t.Cell(13, 1).Range.Text = "82-13="
$t.Cell(13, 2).Range.Text = "59+9="
$t.Cell(13, 3).Range.Text = "85-69="
$t.Cell(13, 4).Range.Text = "6+39="
$t.Cell(13, 5).Range.Text = "96-17="
$t.Cell(14, 1).Range.Text = "94-87="
$t.Cell(14, 2).Range.Text = "64-5="
$t.Cell(14, 3).Range.Text = "8+17="
$t.Cell(14, 4).Range.Text = "49+49="
$t.Cell(14, 5).Range.Text = "92-53="
$t.Cell(15, 1).Range.Text = "73+8="
$t.Cell(15, 2).Range.Text = "28+18="
$t.Cell(15, 3).Range.Text = "37+56="
$t.Cell(15, 4).Range.Text = "83-9="
$t.Cell(15, 5).Range.Text = "47+16="
$t.Cell(16, 1).Range.Text = "80-27="
$t.Cell(16, 2).Range.Text = "82-47="
$t.Cell(16, 3).Range.Text = "24-9="
$t.Cell(16, 4).Range.Text = "55+37="
$t.Cell(16, 5).Range.Text = "53-49="
$t.Cell(17, 1).Range.Text = "71-26="
$t.Cell(17, 2).Range.Text = "48+38="
$t.Cell(17, 3).Range.Text = "58+3="
$t.Cell(17, 4).Range.Text = "62-58="
$t.Cell(17, 5).Range.Text = "52-18="
$t.Cell(18, 1).Range.Text = "90-36="
$t.Cell(18, 2).Range.Text = "60-21="
$t.Cell(18, 3).Range.Text = "91-17="
$t.Cell(18, 4).Range.Text = "85-26="
$t.Cell(18, 5).Range.Text = "72-3="
$t.Cell(19, 1).Range.Text = "9+3="
$t.Cell(19, 2).Range.Text = "85-47="
$t.Cell(19, 3).Range.Text = "75+8="
$t.Cell(19, 4).Range.Text = "25-18="
$t.Cell(19, 5).Range.Text = "56+9="
$t.Cell(20, 1).Range.Text = "4+59="
$t.Cell(20, 2).Range.Text = "51-24="
$t.Cell(20, 3).Range.Text = "62-7="
$t.Cell(20, 4).Range.Text = "17+59="
$t.Cell(20, 5).Range.Text = "6+88="

Write-Host "Updated date + $($t.Rows.Count * $t.Columns.Count) table cells"
